# Update odds data for row 3 (Liverpool M. - Wanderers) and row 7
# (Cardiff Metropolitan - Penybont) on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 3 updates ----
$ws.Range("G3").Value  = 2.2
$ws.Range("I3").Value  = 3.1
$ws.Range("K3").Value  = 2.2
$ws.Range("L3").Value  = 3.6
$ws.Range("M3").Value  = 1.06
$ws.Range("N3").Value  = 10
$ws.Range("U3").Value  = 1.73
$ws.Range("V3").Value  = 2
$ws.Range("W3").Value  = 8
$ws.Range("X3").Value  = 11
$ws.Range("Y3").Value  = 9.5
$ws.Range("Z3").Value  = 21
$ws.Range("AB3").Value = 26
$ws.Range("AC3").Value = 10
$ws.Range("AE3").Value = 13
$ws.Range("AF3").Value = 41
$ws.Range("AH3").Value = 15
$ws.Range("AI3").Value = 11
$ws.Range("AK3").Value = 23
$ws.Range("AM3").Value = 201
$ws.Range("AP3").Value = 21
$ws.Range("AX3").Value = 17
$ws.Range("BB3").Value = 151

# ---- Row 7 updates ----
$ws.Range("G7").Value  = 4.6
$ws.Range("H7").Value  = 3.6
$ws.Range("J7").Value  = 4.9
$ws.Range("K7").Value  = 2.15
$ws.Range("L7").Value  = 2.25
$ws.Range("N7").Value  = 7.3
$ws.Range("O7").Value  = 1.3
$ws.Range("P7").Value  = 3.2
$ws.Range("Q7").Value  = 1.91
$ws.Range("R7").Value  = 1.83
$ws.Range("S7").Value  = 1.4
$ws.Range("T7").Value  = 2.72
$ws.Range("U7").Value  = 1.88
$ws.Range("V7").Value  = 1.83
$ws.Range("Z7").Value  = 80
$ws.Range("AA7").Value = 45
$ws.Range("AB7").Value = 50
$ws.Range("AC7").Value = 7.3
$ws.Range("AD7").Value = 7.1
$ws.Range("AE7").Value = 16.5
$ws.Range("AF7").Value = 80
$ws.Range("AG7").Value = 6.5
$ws.Range("AH7").Value = 7.5
$ws.Range("AL7").Value = 28
$ws.Range("AO7").Value = 27
$ws.Range("AP7").Value = 32
$ws.Range("AS7").Value = 450
$ws.Range("AT7").Value = 2.72
$ws.Range("AU7").Value = 7.7
$ws.Range("AV7").Value = 75
$ws.Range("AW7").Value = 3.5
$ws.Range("AY7").Value = 18.5
$ws.Range("AZ7").Value = 28
$ws.Range("BB7").Value = 250

$wb.Save()
